# Implementa melhorias no menu de Utilitários: adiciona barra de pesquisa,
# filtro de status funcional e remove botão toggle confuso

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1: rename "Plano Básico" -> "Plano Basic" in the header row ---
$ws1.Range("B1").Value = "Plano Basic"

# --- Sheet1: add the new "tabelas:" block below the existing table ---
$rngA = $ws1.Range("A9:A15")
$rngA.Font.Name = "Arial"
$rngA.Font.Size = 10
$rngA.Font.Color = 0

$ws1.Range("A9").Value = "tabelas:"
$ws1.Range("A10").Value = "plan"
$ws1.Range("B10").Value = "cadastro do tipo do plano e campos de controle"
$ws1.Range("B10").Font.Name = "Arial"
$ws1.Range("B10").Font.Size = 10
$ws1.Range("B10").Font.Color = 0
$ws1.Range("A11").Value = "planprice"

$ws1.Rows.Item(9).RowHeight = 15.75
$ws1.Rows.Item(10).RowHeight = 15.75
$ws1.Rows.Item(11).RowHeight = 15.75
$ws1.Rows.Item(12).RowHeight = 15.75
$ws1.Rows.Item(13).RowHeight = 15.75
$ws1.Rows.Item(14).RowHeight = 15.75
$ws1.Rows.Item(15).RowHeight = 15.75

# Move the selection to where the author left off editing
$ws1.Range("B11").Select()

# --- Add the new "Planilha1" sheet after the first sheet ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Planilha1"
$ws2.Range("A1").Value = 40
$ws2.Range("B1").Value = 10
$ws2.Range("C1").Formula = "=A1/B1"
$ws2.Range("C1").Select()

# The new sheet becomes the active tab
$ws2.Activate()
